$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (dates stored as Excel serial numbers, column A uses date format from row above)
# Row 174 -> 2021-09-15 (serial 44454)
$ws.Cells.Item(174, 1).Value = 44454
$ws.Cells.Item(174, 2).Value = 1602737
$ws.Cells.Item(174, 3).Value = 1480519
$ws.Cells.Item(174, 4).Value = 921386
$ws.Cells.Item(174, 5).Value = 700142
$ws.Cells.Item(174, 6).Value = 39750
$ws.Cells.Item(174, 7).Value = 40023

# Row 175 -> 2021-09-16 (serial 44455)
$ws.Cells.Item(175, 1).Value = 44455
$ws.Cells.Item(175, 2).Value = 1602737
$ws.Cells.Item(175, 3).Value = 1489443
$ws.Cells.Item(175, 4).Value = 921386
$ws.Cells.Item(175, 5).Value = 708689
$ws.Cells.Item(175, 6).Value = 39750
$ws.Cells.Item(175, 7).Value = 40023

# Apply the same date number format as the rest of column A to the new cells
$ws.Range("A174:A175").NumberFormat = "yyyy\-mm\-dd"

# Update view: scroll position and active selection, matching the diff
$ws.Application.ActiveWindow.ScrollRow = 156
$ws.Range("D175").Select()
